$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "40÷7="  # was "61÷9="
$t.Cell(1,2).Range.Text = "75÷7="  # was "27÷3="
$t.Cell(1,3).Range.Text = "82÷3="  # was "53÷7="
$t.Cell(1,4).Range.Text = "43÷2="  # was "68÷8="
$t.Cell(1,5).Range.Text = "69÷7="  # was "56÷4="
$t.Cell(5,1).Range.Text = "68÷4="  # was "64÷5="
$t.Cell(5,2).Range.Text = "89÷5="  # was "49÷6="
$t.Cell(5,3).Range.Text = "69÷5="  # was "99÷7="
$t.Cell(5,4).Range.Text = "21÷9="  # was "99÷3="
$t.Cell(5,5).Range.Text = "55÷2="  # was "29÷8="
$t.Cell(9,1).Range.Text = "60÷7="  # was "79÷3="
$t.Cell(9,2).Range.Text = "14÷7="  # was "40÷8="
$t.Cell(9,3).Range.Text = "19÷2="  # was "98÷6="
$t.Cell(9,4).Range.Text = "39÷6="  # was "66÷6="
$t.Cell(9,5).Range.Text = "34÷9="  # was "64÷8="
$t.Cell(13,1).Range.Text = "11÷8="  # was "96÷5="
$t.Cell(13,2).Range.Text = "64÷2="  # was "34÷2="
$t.Cell(13,3).Range.Text = "26÷3="  # was "90÷3="
$t.Cell(13,4).Range.Text = "34÷5="  # was "41÷6="
$t.Cell(13,5).Range.Text = "31÷8="  # was "27÷5="
$t.Cell(17,1).Range.Text = "51÷2="  # was "41÷6="
$t.Cell(17,2).Range.Text = "29÷2="  # was "46÷6="
$t.Cell(17,3).Range.Text = "43÷8="  # was "92÷4="
$t.Cell(17,4).Range.Text = "90÷7="  # was "52÷4="
$t.Cell(17,5).Range.Text = "56÷6="  # was "43÷4="
